$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the first new weekly record at row 102 ---------------------
# This pushes the existing rows 102..199 down to 103..200.
$ws.Rows.Item(102).Insert()

$row102 = @(
    8,
    "Terminal La Palmera de La Serena",
    "Coquimbo",
    44554,
    4,
    100112012,
    "Espinaca",
    "Sin especificar",
    "Primera",
    2800,
    400,
    500,
    450,
    "`$/atado 300 a 500 gramos",
    "Provincia del Elquí",
    900,
    0.5,
    "Hortaliza"
)
for ($i = 0; $i -lt $row102.Length; $i++) {
    $ws.Cells.Item(102, $i + 1).Value = $row102[$i]
}

# --- Insert the second new weekly record at row 141 ---------------------
# (position counted AFTER the row-102 insert above) pushing the rows that
# were 141..200 down to 142..201.
$ws.Rows.Item(141).Insert()

$row141 = @(
    8,
    "Terminal La Palmera de La Serena",
    "Coquimbo",
    44553,
    4,
    100112012,
    "Espinaca",
    "Sin especificar",
    "Primera",
    3000,
    400,
    500,
    450,
    "`$/atado 300 a 500 gramos",
    "Provincia del Elquí",
    900,
    0.5,
    "Hortaliza"
)
for ($i = 0; $i -lt $row141.Length; $i++) {
    $ws.Cells.Item(141, $i + 1).Value = $row141[$i]
}
